$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan")

# Add the new "File Path" header in column Y (row 1)
$ws.Range("Y1").Value = "File Path"

# Update the selected / visible range to reflect the new column
# (scroll so column J is the leftmost visible column, matching topLeftCell="J1")
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 10
$ws.Range("Y2").Select()

# Adjust column widths for A and B (closest achievable values given the
# engine's internal column-width quantization)
$ws.Columns.Item(1).ColumnWidth = 9.83
$ws.Columns.Item(2).ColumnWidth = 21.67
